# "added class one deliverables"
# Adds a new "Class 1" style pair of deliverable rows (Flow Diagram / Functionality
# Definition.txt) with assignment, start/due dates and a note, on the
# "Hardware Development Process" sheet, renaming the two rows that used to read
# "Feature Block Diagram" / "Product Definition.txt". Also switches the active
# sheet/selection from "Professionalism" to "Hardware Development Process".

$wb = $excel.ActiveWorkbook
$wsProf = $wb.Worksheets.Item("Professionalism")
$wsHw = $wb.Worksheets.Item("Hardware Development Process")

# --- Row 6: Flow Diagram -------------------------------------------------
$wsHw.Range("A6").Value = "Flow Diagram"
$wsHw.Range("B6").Value = 1
$wsHw.Range("C6").Value = [DateTime]"2017-05-22"
$wsHw.Range("C6").NumberFormat = "d-mmm"
$wsHw.Range("D6").Value = [DateTime]"2017-05-24"
$wsHw.Range("D6").NumberFormat = "d-mmm"
$wsHw.Range("E6").Value = "End of Day"

# --- Row 7: Functionality Definition.txt --------------------------------
$wsHw.Range("A7").Value = "Functionality Definition.txt"
$wsHw.Range("B7").Value = 1
$wsHw.Range("C7").Value = [DateTime]"2017-05-22"
$wsHw.Range("C7").NumberFormat = "d-mmm"
$wsHw.Range("D7").Value = [DateTime]"2017-05-24"
$wsHw.Range("D7").NumberFormat = "d-mmm"
$wsHw.Range("E7").Value = "End of Day"

# --- Switch the active sheet / selection ---------------------------------
# Before: "Professionalism" tab selected, selection on A9.
# After: "Hardware Development Process" tab selected, selection on E8.
$wsHw.Activate()
$wsHw.Range("E8").Select()
$wsProf.Range("A9").Select()
$wsHw.Activate()
